# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" values replacing the old Strike# values in column G (rows 2-14)
$kValues = @(1, 0, 0, 1, 1, 0, 0, 2, 0, 0, 1, 2, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
